$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 32 ("fit time") values for columns B through X (re-run of the fit)
$row32 = @{
    "B32" = 118.601680273
    "C32" = 113.188978118
    "D32" = 20.34941276800004
    "E32" = 46.95743113200001
    "F32" = 29.32962865399998
    "G32" = 329.207562928
    "H32" = 25.31158334899999
    "I32" = 92.96131702000002
    "J32" = 18.88899304299991
    "K32" = 30.78315490500006
    "L32" = 24.80936239800008
    "M32" = 116.4896158510001
    "N32" = 36.24019994700006
    "O32" = 28.27962081400005
    "P32" = 25.301190787
    "Q32" = 50.25439808900001
    "R32" = 27.21574131800003
    "S32" = 92.84372331300005
    "T32" = 32.17093922499998
    "U32" = 60.62963718900005
    "V32" = 23.45565839100004
    "W32" = 55.7932924459999
    "X32" = 27.39885405900009
}
foreach ($addr in $row32.Keys) {
    $ws.Range($addr).Value = $row32[$addr]
}

# Add new column Y values (new data-translator run results)
$yValues = @{
    "Y10" = 0
    "Y11" = 1
    "Y12" = 0.9792250506952584
    "Y13" = 1.196396870043869
    "Y14" = -0.000009830137559319003
    "Y15" = -0.000009830137559319003
    "Y16" = -0.2798935818792961
    "Y17" = 1.878994870569232
    "Y18" = 3.521801081936098
    "Y31" = 380
    "Y32" = 72.48427583800003
    "Y33" = 12.53529613792121
    "Y34" = 3
    "Y35" = 4.178432045973736
    "Y37" = 189
    "Y38" = 91
    "Y39" = 98
    "Y40" = 0.3061224489795918
    "Y41" = 0.07142857142857142
    "Y42" = 0.9896679335900602
    "Y43" = 0.01564133576666851
}
foreach ($addr in $yValues.Keys) {
    $ws.Range($addr).Value = $yValues[$addr]
}
